$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "B2"  = 6.26
    "A3"  = -21.932
    "E3"  = 16.325
    "E12" = 17.889
    "A14" = -21.559
    "A21" = -20.24
    "A23" = -20.536
    "E24" = 17.183
    "A25" = -21.217
    "B25" = 6.106
    "A26" = -21.006
    "B27" = 5.835
    "A29" = -21.219
    "B31" = 5.918
    "B39" = 7.179
    "B48" = 5.274
    "E50" = 16.474
    "B51" = 6.236
    "B52" = 5.399
    "A53" = -21.912
    "E53" = 16.691
    "B55" = 4.86
    "B56" = 4.977
    "A57" = -22.053
    "B57" = 5.415999999999999
    "E57" = 16.416
    "A59" = -22.5
    "E61" = 16.67
    "E63" = 17.621
    "A69" = -21.462
    "E70" = 17.744
    "B73" = 6.804
    "A79" = -21.129
    "A83" = -22.006
    "E86" = 16.394
    "B89" = 5.500999999999999
    "B90" = 5.767
    "A91" = -21.551
    "B92" = 5.852
    "A93" = -21.379
    "E98" = 16.346
    "E100" = 16.603
    "E102" = 16.446
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
